$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "Price"/"Volume(1h)" columns (D/E) were refreshed by the scraper.
# Values are plain strings in the source file (e.g. "1.003", "28.538.40",
# "  -2.31%  "), so writes go through a helper cell pre-formatted as Text (@)
# and pasted in as values-only - this keeps Excel from reinterpreting strings
# like "1.003" as numbers while leaving the destination cell style untouched.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

function Set-TextValue([string]$addr, [string]$text) {
    $helper.Value = $text
    $helper.Copy()
    $ws.Range($addr).PasteSpecial(-4163) # xlPasteValues
}

Set-TextValue "D2" "28.538.40"
Set-TextValue "E2" "  -2.31%  "
Set-TextValue "D3" "1.793.25"
Set-TextValue "E3" "  -1.86%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  -0.27%  "
Set-TextValue "D5" "231.40"
Set-TextValue "E5" "  -1.19%  "
Set-TextValue "E6" "  -0.85%  "
Set-TextValue "D7" "1.004"
Set-TextValue "E7" "  -0.17%  "
Set-TextValue "D8" "0.2772"
Set-TextValue "E8" "  +1.22%  "
Set-TextValue "D9" "23.40"
Set-TextValue "E9" "  +1.26%  "
Set-TextValue "D10" "0.06754"
Set-TextValue "E10" "  -2.72%  "
Set-TextValue "D11" "0.07558"
Set-TextValue "E11" "  -0.80%  "
Set-TextValue "D12" "1.793.54"
Set-TextValue "E12" "  -2.30%  "
Set-TextValue "D13" "4.788"
Set-TextValue "E13" "  +0.92%  "
Set-TextValue "D14" "0.6115"
Set-TextValue "E14" "  -1.55%  "
Set-TextValue "D15" "2.036.51"
Set-TextValue "E15" "  -1.84%  "
Set-TextValue "D16" "75.70"
Set-TextValue "E16" "  -3.15%  "
Set-TextValue "D17" "0.000008875"
Set-TextValue "E17" "  -8.28%  "
Set-TextValue "D18" "28.534.93"
Set-TextValue "E18" "  -1.37%  "
Set-TextValue "D19" "5.418"
Set-TextValue "E19" "  -4.72%  "
Set-TextValue "E20" "  -0.18%  "
Set-TextValue "D21" "208.63"
Set-TextValue "E21" "  -5.53%  "
Set-TextValue "D22" "11.46"
Set-TextValue "E22" "  -0.38%  "
Set-TextValue "D23" "6.822"
Set-TextValue "E23" "  -0.22%  "
Set-TextValue "E24" "  -0.18%  "
Set-TextValue "D25" "152.61"
Set-TextValue "E25" "  -1.96%  "
Set-TextValue "D26" "7.997"
Set-TextValue "E26" "  +1.13%  "
Set-TextValue "D27" "0.1264"
Set-TextValue "E27" "  -1.68%  "
Set-TextValue "D28" "16.43"
Set-TextValue "E28" "  -0.07%  "
Set-TextValue "D29" "1.417"
Set-TextValue "E29" "  -1.36%  "
Set-TextValue "D30" "0.06119"
Set-TextValue "E30" "  -7.92%  "
Set-TextValue "E31" "  -1.16%  "
Set-TextValue "D32" "3.786"
Set-TextValue "E32" "  -0.85%  "
Set-TextValue "D33" "3.760"
Set-TextValue "E33" "  +0.55%  "
Set-TextValue "D34" "1.726"
Set-TextValue "E34" "  +1.33%  "
Set-TextValue "D35" "1.049"
Set-TextValue "E35" "  -3.62%  "
Set-TextValue "D36" "0.6418"
Set-TextValue "E36" "  +0.52%  "
Set-TextValue "D37" "2.503"
Set-TextValue "E37" "  -1.64%  "
Set-TextValue "D38" "2.705"
Set-TextValue "E38" "  -1.14%  "
Set-TextValue "D39" "0.01689"
Set-TextValue "E39" "  -2.11%  "
Set-TextValue "D40" "1.149.43"
Set-TextValue "E40" "  -2.79%  "
Set-TextValue "D41" "6.318"
Set-TextValue "E41" "  -2.78%  "
Set-TextValue "D42" "0.8740"
Set-TextValue "E42" "  -2.99%  "
Set-TextValue "E43" "  -0.17%  "
Set-TextValue "D44" "100.48"
Set-TextValue "E44" "  +0.31%  "
Set-TextValue "D45" "1.947.02"
Set-TextValue "E45" "  -1.65%  "
Set-TextValue "D46" "60.20"
Set-TextValue "E46" "  -2.73%  "
Set-TextValue "E47" "  -3.55%  "
Set-TextValue "D48" "1.586"
Set-TextValue "E48" "  +1.47%  "
Set-TextValue "D49" "8.371"
Set-TextValue "E49" "  -0.86%  "
Set-TextValue "D50" "0.05450"
Set-TextValue "E50" "  -1.93%  "
Set-TextValue "D51" "0.4472"
Set-TextValue "E51" "  -1.84%  "

$helper.Clear()
$excel.CutCopyMode = $false
